# Auto-generated: refresh market-price derived columns (H-N) across all job sheets
# Values sourced from latest scheduled market data pull.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @{
    "H11" = 439
    "I11" = 439
    "K11" = 439
    "M11" = -299
    "H19" = 2379.8
    "J19" = 499
    "L19" = 499
    "N19" = -849
    "I28" = 620.6667
    "J28" = 600
    "K28" = 620.6667
    "L28" = 600
    "M28" = -135.6667
    "N28" = -1570
    "H87" = 5000
    "J87" = 0
    "L87" = 0
    "H90" = 5000
    "J90" = 0
    "L90" = 0
    "H94" = 10832.333
    "I94" = 10832.333
    "K94" = 10832.333
    "M94" = -10381.333
    "H112" = 3112.3157
    "J112" = 3207.5557
    "L112" = 9622.667099999999
    "N112" = -11838.6671
    "H116" = 2233.3333
    "I116" = 2233.3333
    "J116" = 0
    "K116" = 2233.3333
    "L116" = 0
    "M116" = 1208.6667
    "H137" = 1381.2727
    "I137" = 1381.2727
    "J137" = 0
    "K137" = 4143.8181
    "L137" = 0
    "M137" = -1593.8181
    "H138" = 4212.4
    "J138" = 3958.2222
    "L138" = 11874.6666
    "N138" = -22154.6666
}
foreach ($addr in $ALC_updates.Keys) {
    $ws.Range($addr).Value = $ALC_updates[$addr]
}
foreach ($addr in @("N87","N90","N116","N137")) {
    $ws.Range($addr).ClearContents()
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @{
    "H31" = 30000
    "I31" = 30000
    "K31" = 30000
    "M31" = -29706
    "H32" = 3595.158
    "I32" = 3595.158
    "K32" = 3595.158
    "M32" = -3308.158
    "H34" = 2500
    "I34" = 2500
    "K34" = 2500
    "M34" = -2229
    "H45" = 4560
    "I45" = 1673
    "K45" = 1673
    "M45" = -1296
    "H97" = 66
    "I97" = 76
    "J97" = 52.666668
    "K97" = 76
    "L97" = 52.666668
    "M97" = 420
    "N97" = -1044.666668
    "H132" = 2465.3333
    "I132" = 3396
    "J132" = 2000
    "K132" = 10188
    "L132" = 6000
    "M132" = -7658
    "N132" = -11060
    "H141" = 70000
    "J141" = 70000
    "L141" = 70000
    "N141" = -80360
}
foreach ($addr in $ARM_updates.Keys) {
    $ws.Range($addr).Value = $ARM_updates[$addr]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @{
    "H123" = 42500
    "I123" = 45000
    "K123" = 45000
    "M123" = -40100
}
foreach ($addr in $BSM_updates.Keys) {
    $ws.Range($addr).Value = $BSM_updates[$addr]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @{
    "H31" = 2423.5334
    "I31" = 1797.2307
    "K31" = 1797.2307
    "M31" = -1502.2307
    "H34" = 2423.5334
    "I34" = 1797.2307
    "K34" = 1797.2307
    "M34" = -1595.2307
    "H58" = 1782.7142
    "I58" = 1549.75
    "J58" = 2093.3333
    "K58" = 1549.75
    "L58" = 2093.3333
    "M58" = -1346.75
    "N58" = -2499.3333
    "H122" = 3010
    "I122" = 3010
    "K122" = 9030
    "M122" = -6580
    "H132" = 4839.353
    "I132" = 5016.8125
    "J132" = 2000
    "K132" = 15050.4375
    "L132" = 6000
    "M132" = -12520.4375
    "N132" = -11060
    "H136" = 1782.7142
    "I136" = 1549.75
    "J136" = 2093.3333
    "K136" = 4649.25
    "L136" = 6279.999899999999
    "M136" = -2099.25
    "N136" = -11379.9999
}
foreach ($addr in $CRP_updates.Keys) {
    $ws.Range($addr).Value = $CRP_updates[$addr]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @{
    "H12" = 1078.3846
    "I12" = 105.666664
    "J12" = 1912.1428
    "K12" = 316.999992
    "L12" = 5736.428400000001
    "M12" = -143.999992
    "N12" = -6082.428400000001
    "H81" = 15000
    "J81" = 15000
    "L81" = 45000
    "N81" = -47246
    "H84" = 15000
    "J84" = 15000
    "L84" = 135000
    "N84" = -146232
    "H86" = 861.1818
    "I86" = 539
    "J86" = 1425
    "K86" = 1617
    "L86" = 4275
    "M86" = -431
    "N86" = -6647
    "H89" = 861.1818
    "I89" = 539
    "J89" = 1425
    "K89" = 4851
    "L89" = 12825
    "M89" = 1077
    "N89" = -24681
    "H99" = 8175
    "I99" = 6262.5
    "K99" = 18787.5
    "M99" = -16541.5
}
foreach ($addr in $CUL_updates.Keys) {
    $ws.Range($addr).Value = $CUL_updates[$addr]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @{
    "H63" = 50000
    "J63" = 50000
    "L63" = 50000
    "N63" = -51372
    "H66" = 50000
    "J66" = 50000
    "L66" = 150000
    "N66" = -156864
    "H104" = 0
    "J104" = 0
    "L104" = 0
    "H132" = 4989
    "I132" = 0
    "K132" = 0
    "H136" = 24730.2
    "J136" = 24730.2
    "L136" = 74190.60000000001
    "N136" = -79290.60000000001
}
foreach ($addr in $GSM_updates.Keys) {
    $ws.Range($addr).Value = $GSM_updates[$addr]
}
foreach ($addr in @("N104","M132")) {
    $ws.Range($addr).ClearContents()
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @{
    "H22" = 2134.818
    "I22" = 1283.2858
    "K22" = 1283.2858
    "M22" = -988.2858000000001
    "H27" = 2134.818
    "I27" = 1283.2858
    "K27" = 1283.2858
    "M27" = -1176.2858
    "H40" = 1300
    "I40" = 1300
    "K40" = 1300
    "M40" = -1164
    "H46" = 1983.5927
    "I46" = 1174.4117
    "J46" = 3359.2
    "K46" = 1174.4117
    "L46" = 3359.2
    "M46" = -986.4117000000001
    "N46" = -3735.2
    "H61" = 2029.5714
    "I61" = 1970
    "K61" = 1970
    "M61" = -1768
    "H101" = 29993.25
    "J101" = 29993.25
    "L101" = 29993.25
    "N101" = -36483.25
    "H110" = 30000
    "J110" = 30000
    "L110" = 30000
    "N110" = -38180
    "H113" = 2029.5714
    "I113" = 1970
    "K113" = 1970
    "M113" = 200
    "H122" = 4563.5654
    "I122" = 4429.6313
    "J122" = 5199.75
    "K122" = 13288.8939
    "L122" = 15599.25
    "M122" = -10838.8939
    "N122" = -20499.25
}
foreach ($addr in $LTW_updates.Keys) {
    $ws.Range($addr).Value = $LTW_updates[$addr]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @{
    "H39" = 0
    "J39" = 0
    "L39" = 0
    "H122" = 0
    "I122" = 0
    "K122" = 0
    "H136" = 1121.4286
    "I136" = 1180.0555
    "J136" = 769.6667
    "K136" = 3540.1665
    "L136" = 2309.0001
    "M136" = -990.1664999999998
    "N136" = -7409.0001
}
foreach ($addr in $WVR_updates.Keys) {
    $ws.Range($addr).Value = $WVR_updates[$addr]
}
foreach ($addr in @("N39","M122")) {
    $ws.Range($addr).ClearContents()
}
